$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add P1 and Q1 header cells, copying O1 formatting (bold, border, centered) then setting values
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Row 2
$ws.Range("B2").Value = 24.57726700690615
$ws.Range("C2").Value = 18.72994801111324
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 27.81763776297332
$ws.Range("F2").Value = 34.59345128714452
$ws.Range("G2").Value = 12.60499037337381
$ws.Range("H2").Value = 2.226386170093845
$ws.Range("I2").Value = 3.304557546145663
$ws.Range("J2").Value = 6.695350876082476
$ws.Range("K2").Value = 9.623843916360043
$ws.Range("L2").Value = 5.047621911049819
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 9.71678259546299

# Row 3
$ws.Range("B3").Value = 22.99890156827629
$ws.Range("C3").Value = 17.65463508543564
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 26.27896525384762
$ws.Range("F3").Value = 33.03845185993406
$ws.Range("G3").Value = 12.57777011490236
$ws.Range("H3").Value = 1.965023899609454
$ws.Range("I3").Value = 3.163326210099665
$ws.Range("J3").Value = 6.809903810025647
$ws.Range("K3").Value = 9.962468897681594
$ws.Range("L3").Value = 4.982860707475792
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 9.87259881446023

# Row 4
$ws.Range("B4").Value = 21.97148992022231
$ws.Range("C4").Value = 16.96454721200382
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 25.28918032972688
$ws.Range("F4").Value = 32.05467966720407
$ws.Range("G4").Value = 12.58564807157257
$ws.Range("H4").Value = 1.798794164275098
$ws.Range("I4").Value = 3.074590549696838
$ws.Range("J4").Value = 6.883792019045265
$ws.Range("K4").Value = 10.17557078370752
$ws.Range("L4").Value = 4.94150011795174
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 9.977568094454355

# Row 5
$ws.Range("B5").Value = 21.53800713154276
$ws.Range("C5").Value = 16.68923968624866
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 24.87957174810676
$ws.Range("F5").Value = 31.64211622732023
$ws.Range("G5").Value = 12.58084871972859
$ws.Range("H5").Value = 1.72951746808881
$ws.Range("I5").Value = 3.038364056757102
$ws.Range("J5").Value = 6.912916832850343
$ws.Range("K5").Value = 10.26310943493849
$ws.Range("L5").Value = 4.923840474996704
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 10.0175665757161

# Row 6
$ws.Range("B6").Value = 21.46496410586501
$ws.Range("C6").Value = 16.65933035065913
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 24.81682540877535
$ws.Range("F6").Value = 31.567300466778
$ws.Range("G6").Value = 12.56351728873148
$ws.Range("H6").Value = 1.71784293422426
$ws.Range("I6").Value = 3.033053177522037
$ws.Range("J6").Value = 6.915546966113953
$ws.Range("K6").Value = 10.27698693041963
$ws.Range("L6").Value = 4.920397835130827
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 10.01826478054767

# Row 7
$ws.Range("B7").Value = 21.96523441356292
$ws.Range("C7").Value = 17.00435178610292
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 25.29956828449076
$ws.Range("F7").Value = 32.03327267972501
$ws.Range("G7").Value = 12.53948213416611
$ws.Range("H7").Value = 1.797649344327751
$ws.Range("I7").Value = 3.075821335989203
$ws.Range("J7").Value = 6.878008056307459
$ws.Range("K7").Value = 10.1748183997973
$ws.Range("L7").Value = 4.939950290871159
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 9.96146559371514

# Row 8
$ws.Range("B8").Value = 24.04459615429953
$ws.Range("C8").Value = 18.41922803896933
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 27.31612404303874
$ws.Range("F8").Value = 34.04378110368916
$ws.Range("G8").Value = 12.5296575039375
$ws.Range("H8").Value = 2.137191738990107
$ws.Range("I8").Value = 3.258185026009045
$ws.Range("J8").Value = 6.725881588280995
$ws.Range("K8").Value = 9.737498543437786
$ws.Range("L8").Value = 5.023936150057978
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 9.746272760409031

# Row 9
$ws.Range("B9").Value = 27.66117815882089
$ws.Range("C9").Value = 20.88205482689624
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 30.88069058968735
$ws.Range("F9").Value = 37.76573017227498
$ws.Range("G9").Value = 12.8037315798755
$ws.Range("H9").Value = 2.758806822897629
$ws.Range("I9").Value = 3.59818214454064
$ws.Range("J9").Value = 6.468455884349515
$ws.Range("K9").Value = 8.922634502459509
$ws.Range("L9").Value = 5.178497872995813
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 9.435396625045497

# Row 10
$ws.Range("B10").Value = 30.03676649040536
$ws.Range("C10").Value = 22.5404816994571
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 32.83320038853212
$ws.Range("F10").Value = 39.92670302153588
$ws.Range("G10").Value = 13.00540150531295
$ws.Range("H10").Value = 3.165338493054085
$ws.Range("I10").Value = 3.83703046154572
$ws.Range("J10").Value = 6.281575949554677
$ws.Range("K10").Value = 8.35192387603272
$ws.Range("L10").Value = 5.257248093177494
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 9.213181572293557

# Row 11
$ws.Range("B11").Value = 31.04253063587614
$ws.Range("C11").Value = 23.16560309603797
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 29.9045329237995
$ws.Range("F11").Value = 37.43252330659094
$ws.Range("G11").Value = 11.97299336555735
$ws.Range("H11").Value = 3.738557601642788
$ws.Range("I11").Value = 3.91930384750792
$ws.Range("J11").Value = 6.058489210292378
$ws.Range("K11").Value = 8.164901783918628
$ws.Range("L11").Value = 5.136663981206363
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 8.703304439937014

# Row 12
$ws.Range("B12").Value = 31.40952946593408
$ws.Range("C12").Value = 23.30019614759523
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 27.00365398019699
$ws.Range("F12").Value = 34.86962019202082
$ws.Range("G12").Value = 11.0979190284719
$ws.Range("H12").Value = 4.73911373499895
$ws.Range("I12").Value = 3.937554566923497
$ws.Range("J12").Value = 5.917957640235497
$ws.Range("K12").Value = 8.14900887504684
$ws.Range("L12").Value = 5.132449140282676
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 8.340056300937608

# Row 13
$ws.Range("B13").Value = 31.31953224133494
$ws.Range("C13").Value = 23.12236430613715
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 23.89490671930609
$ws.Range("F13").Value = 32.01775985187801
$ws.Range("G13").Value = 10.20657711058226
$ws.Range("H13").Value = 5.895021633247699
$ws.Range("I13").Value = 3.910685635960187
$ws.Range("J13").Value = 5.823857976844272
$ws.Range("K13").Value = 8.25612170359543
$ws.Range("L13").Value = 5.212575323280736
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 8.044012846967886

# Row 14
$ws.Range("B14").Value = 31.0485398804624
$ws.Range("C14").Value = 22.8644799312799
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 21.60363715531179
$ws.Range("F14").Value = 29.86176544606688
$ws.Range("G14").Value = 9.564468606101846
$ws.Range("H14").Value = 6.761429188414043
$ws.Range("I14").Value = 3.872398408595491
$ws.Range("J14").Value = 5.779017449424177
$ws.Range("K14").Value = 8.387662002657963
$ws.Range("L14").Value = 5.317723268075151
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 7.868096112625475

# Row 15
$ws.Range("B15").Value = 30.88355111870696
$ws.Range("C15").Value = 22.74769308135729
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 20.99908439300258
$ws.Range("F15").Value = 29.26822250980261
$ws.Range("G15").Value = 9.393014060470232
$ws.Range("H15").Value = 6.95902531803341
$ws.Range("I15").Value = 3.854238426408473
$ws.Range("J15").Value = 5.776486048982838
$ws.Range("K15").Value = 8.444687110971639
$ws.Range("L15").Value = 5.345832786256486
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 7.835926021765048

# Row 16
$ws.Range("B16").Value = 29.92918421238997
$ws.Range("C16").Value = 22.11510516920581
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 20.60918212332188
$ws.Range("F16").Value = 28.77014270703127
$ws.Range("G16").Value = 9.427406341968293
$ws.Range("H16").Value = 6.677418368525549
$ws.Range("I16").Value = 3.759581571691426
$ws.Range("J16").Value = 5.874264565697755
$ws.Range("K16").Value = 8.64906480718795
$ws.Range("L16").Value = 5.291246705667406
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 7.988153169432183

# Row 17
$ws.Range("B17").Value = 29.33135163784739
$ws.Range("C17").Value = 21.76512258254013
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 21.57299530166553
$ws.Range("F17").Value = 29.57238374144572
$ws.Range("G17").Value = 9.801012041046008
$ws.Range("H17").Value = 5.948579654227844
$ws.Range("I17").Value = 3.707262333731594
$ws.Range("J17").Value = 5.972040753554424
$ws.Range("K17").Value = 8.739247022447913
$ws.Range("L17").Value = 5.184207595031291
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 8.194791275103091

# Row 18
$ws.Range("B18").Value = 28.98691447815762
$ws.Range("C18").Value = 21.59627066785843
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 23.83804751478451
$ws.Range("F18").Value = 31.60277831872941
$ws.Range("G18").Value = 10.52756998042087
$ws.Range("H18").Value = 4.811101020803068
$ws.Range("I18").Value = 3.683838281158877
$ws.Range("J18").Value = 6.085281718818544
$ws.Range("K18").Value = 8.745608170341113
$ws.Range("L18").Value = 5.070954060520845
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 8.482925311035903

# Row 19
$ws.Range("B19").Value = 28.87671988693585
$ws.Range("C19").Value = 21.65495712369679
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 27.00678720415926
$ws.Range("F19").Value = 34.43085374440559
$ws.Range("G19").Value = 11.4018662579088
$ws.Range("H19").Value = 3.637007288158027
$ws.Range("I19").Value = 3.692941653684203
$ws.Range("J19").Value = 6.194398509885403
$ws.Range("K19").Value = 8.703606856867529
$ws.Range("L19").Value = 5.041934045968941
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 8.792899508056253

# Row 20
$ws.Range("B20").Value = 29.43069337618095
$ws.Range("C20").Value = 22.22272634970733
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 32.35127332636606
$ws.Range("F20").Value = 39.31049897342385
$ws.Range("G20").Value = 12.80061433308661
$ws.Range("H20").Value = 3.057604041755434
$ws.Range("I20").Value = 3.781744297408231
$ws.Range("J20").Value = 6.308894943845408
$ws.Range("K20").Value = 8.499585590512453
$ws.Range("L20").Value = 5.232106107050295
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 9.212046537255342

# Row 21
$ws.Range("B21").Value = 31.16885129658473
$ws.Range("C21").Value = 23.44522909166952
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 34.45772322659813
$ws.Range("F21").Value = 41.54254075516298
$ws.Range("G21").Value = 13.23369276757824
$ws.Range("H21").Value = 3.401117748663243
$ws.Range("I21").Value = 3.966489194186487
$ws.Range("J21").Value = 6.200525402180237
$ws.Range("K21").Value = 8.061923917991493
$ws.Range("L21").Value = 5.33325851019074
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 9.145380098937798

# Row 22
$ws.Range("B22").Value = 32.25422332983276
$ws.Range("C22").Value = 24.16615423239478
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 35.54754024790277
$ws.Range("F22").Value = 42.77696428711327
$ws.Range("G22").Value = 13.536125148379
$ws.Range("H22").Value = 3.60839289096509
$ws.Range("I22").Value = 4.081364460833766
$ws.Range("J22").Value = 6.135006227867387
$ws.Range("K22").Value = 7.782272768405021
$ws.Range("L22").Value = 5.385834250684543
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 9.115614081341244

# Row 23
$ws.Range("B23").Value = 31.67982969427932
$ws.Range("C23").Value = 23.74529588487158
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 34.955867706839
$ws.Range("F23").Value = 42.13561261830439
$ws.Range("G23").Value = 13.42135245658454
$ws.Range("H23").Value = 3.498291428069433
$ws.Range("I23").Value = 4.017412809989181
$ws.Range("J23").Value = 6.177050197482518
$ws.Range("K23").Value = 7.931145429664142
$ws.Range("L23").Value = 5.359245378122099
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 9.15030664642583

# Row 24
$ws.Range("B24").Value = 29.40338616100395
$ws.Range("C24").Value = 22.1534749652118
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 32.65066012638369
$ws.Range("F24").Value = 39.61656289831325
$ws.Range("G24").Value = 12.97259402700164
$ws.Range("H24").Value = 3.072764164664619
$ws.Range("I24").Value = 3.776685685086167
$ws.Range("J24").Value = 6.334155689249582
$ws.Range("K24").Value = 8.503028369487398
$ws.Range("L24").Value = 5.254123511348564
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 9.279997729872697

# Row 25
$ws.Range("B25").Value = 26.73144120143305
$ws.Range("C25").Value = 20.3093410912366
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 29.98052678241577
$ws.Range("F25").Value = 36.76245902161312
$ws.Range("G25").Value = 12.62378738892559
$ws.Range("H25").Value = 2.59495521707119
$ws.Range("I25").Value = 3.511208782368684
$ws.Range("J25").Value = 6.523352523779099
$ws.Range("K25").Value = 9.137259344308193
$ws.Range("L25").Value = 5.135757549898421
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 9.479094358384536
